$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44322

# Row 3
$ws.Range("D3").Value = 44313
$ws.Range("M3").Value = 120

# Row 4
$ws.Range("D4").Value = 44309
$ws.Range("Q4").Value = '$/caja 14 kilos granel'
$ws.Range("S4").Value = 821
$ws.Range("T4").Value = 14

# Row 5
$ws.Range("D5").Value = 44327
$ws.Range("M5").Value = 60

# Row 6
$ws.Range("D6").Value = 44302
$ws.Range("M6").Value = 80

# Row 7
$ws.Range("D7").Value = 44330
$ws.Range("M7").Value = 60

# Row 9
$ws.Range("D9").Value = 44323
$ws.Range("M9").Value = 80

# Row 10
$ws.Range("D10").Value = 44306
$ws.Range("Q10").Value = '$/caja 10 kilos empedrada'
$ws.Range("S10").Value = 11500
$ws.Range("T10").Value = 1
